$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("glory")

# --- Add new game 12 row of data ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 21
$ws.Range("C12").Value = 13
$ws.Range("D12").Value = 9
$ws.Range("E12").Value = 21
$ws.Range("F12").Value = 14
$ws.Range("G12").Value = 9

$ws.Range("H12").Formula = "=SUM(C`$2:C12)"
$ws.Range("I12").Formula = "=SUM(D`$2:D12)"
$ws.Range("J12").Formula = "=SUM(E`$2:E12)"
$ws.Range("K12").Formula = "=SUM(F`$2:F12)"
$ws.Range("L12").Formula = "=SUM(G`$2:G12)"

$ws.Range("M12").Formula = "=H12-MAX(H12:L12)"
$ws.Range("N12").Formula = "=I12-MAX(H12:L12)"
$ws.Range("O12").Formula = "=J12-MAX(H12:L12)"
$ws.Range("P12").Formula = "=K12-MAX(H12:L12)"
$ws.Range("Q12").Formula = "=L12-MAX(H12:L12)"

$ws.Range("R12").Formula = "=RANK(C12,`$C12:`$G12)"
$ws.Range("S12").Formula = "=RANK(D12,`$C12:`$G12)"
$ws.Range("T12").Formula = "=RANK(E12,`$C12:`$G12)"
$ws.Range("U12").Formula = "=RANK(F12,`$C12:`$G12)"
$ws.Range("V12").Formula = "=RANK(G12,`$C12:`$G12)"

$ws.Range("W12").Formula = "=C12/`$B12"
$ws.Range("X12").Formula = "=D12/`$B12"
$ws.Range("Y12").Formula = "=E12/`$B12"
$ws.Range("Z12").Formula = "=F12/`$B12"
$ws.Range("AA12").Formula = "=G12/`$B12"

$ws.Range("AB12").Formula = "=SUM(C12:G12)"

# --- Update the active sheet / selection to reflect post-edit state ---
# The "glory" sheet becomes the active tab, with the cursor resting on
# A13 (the row just after the newly entered data).
$ws.Activate() | Out-Null
$ws.Range("A13").Select() | Out-Null
